$p = $ppt.ActivePresentation

# -----------------------------------------------------------------
# 1) Table on slide 5 ("B1- TYPES OF FINANCIAL DOCUMENTS") switches
#    from the custom "Table_0" style to the built-in table style
#    {A4792E67-11A1-485D-8B52-521614E600AE}.
# -----------------------------------------------------------------
$slide5 = $p.Slides.Item(5)
$tableShape = $slide5.Shapes.Item(2)
$table = $tableShape.Table
$table.ApplyStyle("{A4792E67-11A1-485D-8B52-521614E600AE}")

# -----------------------------------------------------------------
# 2) The deck's theme colours are swapped: the active theme (the one
#    backing the slide master, currently the "Integral" / "Red
#    Violet" palette) becomes the plain "Office" palette.
#    Theme colour order exposed via ThemeColorScheme:
#      1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink
# -----------------------------------------------------------------
$officeHex = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

$themeColors = $slide5.ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $hex = $officeHex[$i - 1]
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    $rgbVal = $r + ($g * 256) + ($b * 65536)
    $themeColors.Item($i).RGB = $rgbVal
}
